# Resident Services_Requirements.xlsx - add Clarification Description research
# notes to the "Details" sheet (column O) for rows 3-11, matching the
# upstream commit that populated the previously-empty
# "Clarification Description" table column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# --- Row 3 (S.No. 1 - Lock/Unlock UIN for each Auth Type) ---------------
$ws.Range("O3").Value = "User is travelling out of country, personal choice, doesn,t want misuse. Can be anything. It's a user choice. That is the use case. `nCheck if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 4 (S.No. 2 - Request re-print of UIN) ---------------------------
$ws.Range("O4").Value = "User likes multiple copies, use lost old copy. Can be anything. Provision to keep track on no of reprint required, country can reject free request after a X upper limit, X can be zero, country to allow payment based reprint, thus integration point with payment gateway will be required for SI customisation etc. `nvalidation and the Interface for Payment gateway will be part of Resident portal(business Logic)- by SI`nThe service should cater to any service received for Re-printing`nCheck if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 5 (S.No. 3 - Initiate UIN Update) --------------------------------
$ws.Range("O5").Value = "Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 6 (S.No. 4 - Initiate UIN Update - Address) ----------------------
$ws.Range("O6").Value = "Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 7 (S.No. 5 - Track Status of UIN Update) -------------------------
$ws.Range("O7").Value = "Can update any demo graphic information, validations should be done by Portal.`nShould cater to updation of demographic details as done by Reg. client`nCheck if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 8 (S.No. 6 - Track status of UIN Generation) ---------------------
$ws.Range("O8").Value = "think and come back on what all status can the request have and we can rationalise. `nCheck if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 9 (S.No. 7 - Retrieve Lost UIN) -----------------------------------
$ws.Range("O9").Value = "think and come back on what all status can the request have and we can rationalise. `nCheck if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 10 (S.No. 8 - Download UIN) ---------------------------------------
$ws.Range("O10").Value = "YES. Correct understanding. `nCheck if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# --- Row 11 (S.No. 9 - View History of Authentication Requests) -----------
$ws.Range("O11").Value = "Check if OTP authentication is allowed or not- (if a Mobile Numberis associated with multiple UIN/VID then the country can term it as Weak based on the number)"

# Row heights grow/shrink to fit the newly-added text (matches the saved
# workbook's recorded row heights for the Details sheet).
$ws.Rows.Item(3).RowHeight = 246.5
$ws.Rows.Item(4).RowHeight = 319
$ws.Rows.Item(5).RowHeight = 174
$ws.Rows.Item(6).RowHeight = 174
$ws.Rows.Item(8).RowHeight = 246.5
$ws.Rows.Item(9).RowHeight = 174
$ws.Rows.Item(10).RowHeight = 145
$ws.Rows.Item(11).RowHeight = 188.5
$ws.Rows.Item(13).RowHeight = 188.5
$ws.Rows.Item(14).RowHeight = 43.5

# View state: the editor zoomed in on the Details sheet and left the
# selection on the last-edited cell (O11).
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$ws.Range("O11").Select()
